$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A81:G81").Copy($ws.Range("A83:G83"))
$ws.Range("A83").Value = "Year 4"
$ws.Range("B83").Value = "B1D"
$ws.Range("C83").Value = "clinical"
$ws.Range("D83").Value = "'1"
$ws.Range("E83").Value = "'25/11/2025"
$ws.Range("F83").Value = "10:30:00"
$ws.Range("G83").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A84:G84"))
$ws.Range("A84").Value = "Year 4"
$ws.Range("B84").Value = "B1D"
$ws.Range("C84").Value = "clinical"
$ws.Range("D84").Value = "'2"
$ws.Range("E84").Value = "'26/11/2025"
$ws.Range("F84").Value = "10:30:00"
$ws.Range("G84").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A85:G85"))
$ws.Range("A85").Value = "Year 4"
$ws.Range("B85").Value = "B1D"
$ws.Range("C85").Value = "clinical"
$ws.Range("D85").Value = "'3"
$ws.Range("E85").Value = "'27/11/2025"
$ws.Range("F85").Value = "10:30:00"
$ws.Range("G85").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A86:G86"))
$ws.Range("A86").Value = "Year 4"
$ws.Range("B86").Value = "B1D"
$ws.Range("C86").Value = "clinical"
$ws.Range("D86").Value = "'4"
$ws.Range("E86").Value = "'02/12/2025"
$ws.Range("F86").Value = "10:30:00"
$ws.Range("G86").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A87:G87"))
$ws.Range("A87").Value = "Year 4"
$ws.Range("B87").Value = "B1D"
$ws.Range("C87").Value = "clinical"
$ws.Range("D87").Value = "'5"
$ws.Range("E87").Value = "'03/12/2025"
$ws.Range("F87").Value = "10:30:00"
$ws.Range("G87").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A88:G88"))
$ws.Range("A88").Value = "Year 4"
$ws.Range("B88").Value = "B1D"
$ws.Range("C88").Value = "clinical"
$ws.Range("D88").Value = "'6"
$ws.Range("E88").Value = "'04/12/2025"
$ws.Range("F88").Value = "10:30:00"
$ws.Range("G88").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A89:G89"))
$ws.Range("A89").Value = "Year 4"
$ws.Range("B89").Value = "B1D"
$ws.Range("C89").Value = "clinical"
$ws.Range("D89").Value = "'7"
$ws.Range("E89").Value = "'09/12/2025"
$ws.Range("F89").Value = "10:30:00"
$ws.Range("G89").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A90:G90"))
$ws.Range("A90").Value = "Year 4"
$ws.Range("B90").Value = "B1D"
$ws.Range("C90").Value = "clinical"
$ws.Range("D90").Value = "'8"
$ws.Range("E90").Value = "'10/12/2025"
$ws.Range("F90").Value = "10:30:00"
$ws.Range("G90").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A91:G91"))
$ws.Range("A91").Value = "Year 4"
$ws.Range("B91").Value = "B1D"
$ws.Range("C91").Value = "clinical"
$ws.Range("D91").Value = "'9"
$ws.Range("E91").Value = "'11/12/2025"
$ws.Range("F91").Value = "10:30:00"
$ws.Range("G91").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A92:G92"))
$ws.Range("A92").Value = "Year 4"
$ws.Range("B92").Value = "B1D"
$ws.Range("C92").Value = "clinical"
$ws.Range("D92").Value = "'10"
$ws.Range("E92").Value = "'16/12/2025"
$ws.Range("F92").Value = "10:30:00"
$ws.Range("G92").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A93:G93"))
$ws.Range("A93").Value = "Year 4"
$ws.Range("B93").Value = "B1D"
$ws.Range("C93").Value = "clinical"
$ws.Range("D93").Value = "'11"
$ws.Range("E93").Value = "'17/12/2025"
$ws.Range("F93").Value = "10:30:00"
$ws.Range("G93").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A94:G94"))
$ws.Range("A94").Value = "Year 4"
$ws.Range("B94").Value = "B1D"
$ws.Range("C94").Value = "clinical"
$ws.Range("D94").Value = "'12"
$ws.Range("E94").Value = "'18/12/2025"
$ws.Range("F94").Value = "10:30:00"
$ws.Range("G94").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A95:G95"))
$ws.Range("A95").Value = "Year 4"
$ws.Range("B95").Value = "B1D"
$ws.Range("C95").Value = "clinical"
$ws.Range("D95").Value = "'13"
$ws.Range("E95").Value = "'23/12/2025"
$ws.Range("F95").Value = "10:30:00"
$ws.Range("G95").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A96:G96"))
$ws.Range("A96").Value = "Year 4"
$ws.Range("B96").Value = "B1D"
$ws.Range("C96").Value = "clinical"
$ws.Range("D96").Value = "'14"
$ws.Range("E96").Value = "'24/12/2025"
$ws.Range("F96").Value = "10:30:00"
$ws.Range("G96").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A97:G97"))
$ws.Range("A97").Value = "Year 4"
$ws.Range("B97").Value = "B1D"
$ws.Range("C97").Value = "clinical"
$ws.Range("D97").Value = "'15"
$ws.Range("E97").Value = "'25/12/2025"
$ws.Range("F97").Value = "10:30:00"
$ws.Range("G97").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A98:G98"))
$ws.Range("A98").Value = "Year 4"
$ws.Range("B98").Value = "B1D"
$ws.Range("C98").Value = "clinical"
$ws.Range("D98").Value = "'16"
$ws.Range("E98").Value = "'30/12/2025"
$ws.Range("F98").Value = "10:30:00"
$ws.Range("G98").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A99:G99"))
$ws.Range("A99").Value = "Year 4"
$ws.Range("B99").Value = "B1D"
$ws.Range("C99").Value = "clinical"
$ws.Range("D99").Value = "'17"
$ws.Range("E99").Value = "'31/12/2025"
$ws.Range("F99").Value = "10:30:00"
$ws.Range("G99").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A100:G100"))
$ws.Range("A100").Value = "Year 4"
$ws.Range("B100").Value = "B1D"
$ws.Range("C100").Value = "clinical"
$ws.Range("D100").Value = "'18"
$ws.Range("E100").Value = "'01/01/2026"
$ws.Range("F100").Value = "10:30:00"
$ws.Range("G100").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A101:G101"))
$ws.Range("A101").Value = "Year 4"
$ws.Range("B101").Value = "B1D"
$ws.Range("C101").Value = "clinical"
$ws.Range("D101").Value = "'19"
$ws.Range("E101").Value = "'06/01/2026"
$ws.Range("F101").Value = "10:30:00"
$ws.Range("G101").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A102:G102"))
$ws.Range("A102").Value = "Year 4"
$ws.Range("B102").Value = "B1D"
$ws.Range("C102").Value = "clinical"
$ws.Range("D102").Value = "'20"
$ws.Range("E102").Value = "'07/01/2026"
$ws.Range("F102").Value = "10:30:00"
$ws.Range("G102").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A103:G103"))
$ws.Range("A103").Value = "Year 4"
$ws.Range("B103").Value = "B1D"
$ws.Range("C103").Value = "clinical"
$ws.Range("D103").Value = "'21"
$ws.Range("E103").Value = "'08/01/2026"
$ws.Range("F103").Value = "10:30:00"
$ws.Range("G103").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A104:G104"))
$ws.Range("A104").Value = "Year 4"
$ws.Range("B104").Value = "B1D"
$ws.Range("C104").Value = "clinical"
$ws.Range("D104").Value = "'22"
$ws.Range("E104").Value = "'13/01/2026"
$ws.Range("F104").Value = "10:30:00"
$ws.Range("G104").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A105:G105"))
$ws.Range("A105").Value = "Year 4"
$ws.Range("B105").Value = "B1D"
$ws.Range("C105").Value = "clinical"
$ws.Range("D105").Value = "'23"
$ws.Range("E105").Value = "'14/01/2026"
$ws.Range("F105").Value = "10:30:00"
$ws.Range("G105").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A106:G106"))
$ws.Range("A106").Value = "Year 4"
$ws.Range("B106").Value = "B1D"
$ws.Range("C106").Value = "clinical"
$ws.Range("D106").Value = "'24"
$ws.Range("E106").Value = "'15/01/2026"
$ws.Range("F106").Value = "10:30:00"
$ws.Range("G106").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A107:G107"))
$ws.Range("A107").Value = "Year 4"
$ws.Range("B107").Value = "B1D"
$ws.Range("C107").Value = "clinical"
$ws.Range("D107").Value = "'25"
$ws.Range("E107").Value = "'20/01/2026"
$ws.Range("F107").Value = "10:30:00"
$ws.Range("G107").Value = 180

$ws.Range("A82:G82").Copy($ws.Range("A108:G108"))
$ws.Range("A108").Value = "Year 4"
$ws.Range("B108").Value = "B1D"
$ws.Range("C108").Value = "clinical"
$ws.Range("D108").Value = "'26"
$ws.Range("E108").Value = "'21/01/2026"
$ws.Range("F108").Value = "10:30:00"
$ws.Range("G108").Value = 180

$ws.Range("A81:G81").Copy($ws.Range("A109:G109"))
$ws.Range("A109").Value = "Year 4"
$ws.Range("B109").Value = "B1D"
$ws.Range("C109").Value = "clinical"
$ws.Range("D109").Value = "'27"
$ws.Range("E109").Value = "'22/01/2026"
$ws.Range("F109").Value = "10:30:00"
$ws.Range("G109").Value = 180
